$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.948.70"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "'1.812.98"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'310.25"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.4985"
$ws.Range("E7").Value = "  -2.75%  "
$ws.Range("D8").Value = "'0.3889"
$ws.Range("E8").Value = "  +2.89%  "
$ws.Range("D9").Value = "'0.09689"
$ws.Range("E9").Value = "  +24.72%  "
$ws.Range("D10").Value = "'1.102"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("D11").Value = "'41.02"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "'6.424"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("D13").Value = "'20.50"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Value = "'1.001"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.306"
$ws.Range("E15").Value = "  +1.99%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'1.804.23"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "'0.00001125"
$ws.Range("E17").Value = "  +5.20%  "
$ws.Range("D18").Value = "'92.62"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("D19").Value = "'0.06659"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "'17.12"
$ws.Range("D22").Value = "'5.919"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").Value = "'28.019.81"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "'11.11"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").Value = "'2.240"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'159.02"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "'2.021.90"
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.58"
$ws.Range("E28").Value = "  +2.00%  "
$ws.Range("D29").Value = "'2.392"
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("D30").Value = "'128.08"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("D32").Value = "'1.039"
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "'5.571"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").Value = "'3.627"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").Value = "'0.06719"
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("D36").Value = "'8.985"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").Value = "'0.02330"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").Value = "'0.2136"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").Value = "'4.943"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "'11.24"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "'0.6190"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "'1.147"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "'13.19"
$ws.Range("E44").Value = "  +0.68%  "
$ws.Range("D45").Value = "'0.5872"
$ws.Range("D46").Value = "'3.689"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "'1.279"
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("D48").Value = "'123.05"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("D49").Value = "'1.939"
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").Value = "'1.180"
$ws.Range("E50").Value = "  -2.65%  "
$ws.Range("D51").Value = "'0.06791"
$ws.Range("E51").Value = "  +1.32%  "

Write-Host "Applied 101 cell updates"
